# Applies the "Додаткові критерії вибору" sheet changes:
#  - shift the header picture one column to the left (C -> B)
#  - tighten rows 1-3 to a fixed 13pt height and grow row 4 to 30pt
#  - number the criteria rows (A5:A13) 1..9
#  - drop the stray left-hairline border on C5/C6 so they match the rest
#    of their row
#  - give row 14 a caption in column B

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Додаткові критерії вибору")

# --- move the logo picture from column C to column B (keep same top-left offset) ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 19.5
$shp.Top = 0

# --- row heights ---
$ws.Rows.Item(1).RowHeight = 13
$ws.Rows.Item(2).RowHeight = 13
$ws.Rows.Item(3).RowHeight = 13
$ws.Rows.Item(4).RowHeight = 30

# --- sequence numbers for the criteria rows ---
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 4
$ws.Range("A9").Value = 5
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 7
$ws.Range("A12").Value = 8
$ws.Range("A13").Value = 9

# --- C5/C6 pick up the plain row border instead of the special hairline one ---
$ws.Range("A5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- caption for the results row ---
$ws.Range("B14").Value = "Результати дослідження діяльності обраних компаній за інформацією, отриманою з мережі Інтернет, та додаткової перевірки незалежності компаній"
